$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 462; existing rows 462:495 shift down to 463:496.
$ws.Rows.Item(462).Insert()

# Fill in the new row 462 with the new record's data.
$ws.Range("A462").Value = 6
$ws.Range("B462").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C462").Value = "Metropolitana"
$ws.Range("D462").Value = 44746
$ws.Range("E462").Value = 13
$ws.Range("F462").Value = 100112052
$ws.Range("G462").Value = "Albahaca"
$ws.Range("H462").Value = "Sin especificar"
$ws.Range("I462").Value = "Primera"
$ws.Range("J462").Value = 140
$ws.Range("K462").Value = 4000
$ws.Range("L462").Value = 4500
$ws.Range("M462").Value = 4214
$ws.Range("N462").Value = "$/paquete"
$ws.Range("O462").Value = "Región de Arica y Parinacota"
$ws.Range("P462").Value = 4214
$ws.Range("Q462").Value = 1
$ws.Range("R462").Value = "Hortaliza"

# Match the date cell formatting used by the rest of the D column (custom date format).
$ws.Range("D462").NumberFormat = $ws.Range("D461").NumberFormat
